$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the 3 new header labels for the Brute-force clustering columns (H, I, J) ---
$ws.Range("H1").Value = "ClusterSize(Brute-force)"
$ws.Range("I1").Value = "MinimumPercentIdentity(Brute-force)"
$ws.Range("J1").Value = "Average(Brute-force)"

# --- Copy the header fill style from the existing last header cell (G1) onto H1:J1 ---
[void]$ws.Range("G1").Copy()
[void]$ws.Range("H1:J1").PasteSpecial(-4122)

# Re-assert the header text values (PasteSpecial(Formats) should not disturb them, done defensively)
$ws.Range("H1").Value = "ClusterSize(Brute-force)"
$ws.Range("I1").Value = "MinimumPercentIdentity(Brute-force)"
$ws.Range("J1").Value = "Average(Brute-force)"

# --- Fill in the new Brute-force clustering results for rows 2-51 ---
$ws.Range("H2").Value = 67
$ws.Range("I2").Value = 88
$ws.Range("J2").Value = 93.781999999999996
$ws.Range("H3").Value = 53
$ws.Range("I3").Value = 87
$ws.Range("J3").Value = 95.493499999999997
$ws.Range("H4").Value = 51
$ws.Range("I4").Value = 93
$ws.Range("J4").Value = 97.036900000000003
$ws.Range("H5").Value = 49
$ws.Range("I5").Value = 94
$ws.Range("J5").Value = 96.834999999999994
$ws.Range("H6").Value = 47
$ws.Range("I6").Value = 89
$ws.Range("J6").Value = 96.114699999999999
$ws.Range("H7").Value = 39
$ws.Range("I7").Value = 89
$ws.Range("J7").Value = 95.175399999999996
$ws.Range("H8").Value = 39
$ws.Range("I8").Value = 91
$ws.Range("J8").Value = 97.357600000000005
$ws.Range("H9").Value = 34
$ws.Range("I9").Value = 88
$ws.Range("J9").Value = 95.106999999999999
$ws.Range("H10").Value = 34
$ws.Range("I10").Value = 85
$ws.Range("J10").Value = 95.916200000000003
$ws.Range("H11").Value = 33
$ws.Range("I11").Value = 94
$ws.Range("J11").Value = 96.358000000000004
$ws.Range("H12").Value = 30
$ws.Range("I12").Value = 90
$ws.Range("J12").Value = 94.425299999999993
$ws.Range("H13").Value = 30
$ws.Range("I13").Value = 88
$ws.Range("J13").Value = 92.314899999999994
$ws.Range("H14").Value = 29
$ws.Range("I14").Value = 89
$ws.Range("J14").Value = 94.578800000000001
$ws.Range("H15").Value = 25
$ws.Range("I15").Value = 84
$ws.Range("J15").Value = 91.41
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 90
$ws.Range("J16").Value = 96.036699999999996
$ws.Range("H17").Value = 22
$ws.Range("I17").Value = 94
$ws.Range("J17").Value = 96.29
$ws.Range("H18").Value = 22
$ws.Range("I18").Value = 89
$ws.Range("J18").Value = 93.935100000000006
$ws.Range("H19").Value = 20
$ws.Range("I19").Value = 91
$ws.Range("J19").Value = 95.547399999999996
$ws.Range("H20").Value = 19
$ws.Range("I20").Value = 89
$ws.Range("J20").Value = 94.251499999999993
$ws.Range("H21").Value = 19
$ws.Range("I21").Value = 97
$ws.Range("J21").Value = 97.117000000000004
$ws.Range("H22").Value = 18
$ws.Range("I22").Value = 95
$ws.Range("J22").Value = 97.3399
$ws.Range("H23").Value = 18
$ws.Range("I23").Value = 87
$ws.Range("J23").Value = 92.477099999999993
$ws.Range("H24").Value = 18
$ws.Range("I24").Value = 90
$ws.Range("J24").Value = 95.097999999999999
$ws.Range("H25").Value = 17
$ws.Range("I25").Value = 88
$ws.Range("J25").Value = 94.301500000000004
$ws.Range("H26").Value = 17
$ws.Range("I26").Value = 86
$ws.Range("J26").Value = 91.654399999999995
$ws.Range("H27").Value = 17
$ws.Range("I27").Value = 92
$ws.Range("J27").Value = 95.2059
$ws.Range("H28").Value = 17
$ws.Range("I28").Value = 84
$ws.Range("J28").Value = 90.7059
$ws.Range("H29").Value = 17
$ws.Range("I29").Value = 87
$ws.Range("J29").Value = 92.463200000000001
$ws.Range("H30").Value = 17
$ws.Range("I30").Value = 90
$ws.Range("J30").Value = 95.625
$ws.Range("H31").Value = 16
$ws.Range("I31").Value = 89
$ws.Range("J31").Value = 93.291700000000006
$ws.Range("H32").Value = 16
$ws.Range("I32").Value = 95
$ws.Range("J32").Value = 97.4
$ws.Range("H33").Value = 16
$ws.Range("I33").Value = 88
$ws.Range("J33").Value = 94.291700000000006
$ws.Range("H34").Value = 15
$ws.Range("I34").Value = 90
$ws.Range("J34").Value = 94.628600000000006
$ws.Range("H35").Value = 15
$ws.Range("I35").Value = 88
$ws.Range("J35").Value = 93.6
$ws.Range("H36").Value = 15
$ws.Range("I36").Value = 90
$ws.Range("J36").Value = 95.552400000000006
$ws.Range("H37").Value = 15
$ws.Range("I37").Value = 90
$ws.Range("J37").Value = 97.047600000000003
$ws.Range("H38").Value = 15
$ws.Range("I38").Value = 87
$ws.Range("J38").Value = 94.676199999999994
$ws.Range("H39").Value = 15
$ws.Range("I39").Value = 89
$ws.Range("J39").Value = 93.4
$ws.Range("H40").Value = 15
$ws.Range("I40").Value = 85
$ws.Range("J40").Value = 90.009500000000003
$ws.Range("H41").Value = 14
$ws.Range("I41").Value = 92
$ws.Range("J41").Value = 95.483500000000006
$ws.Range("H42").Value = 14
$ws.Range("I42").Value = 90
$ws.Range("J42").Value = 94.692300000000003
$ws.Range("H43").Value = 14
$ws.Range("I43").Value = 93
$ws.Range("J43").Value = 96.780199999999994
$ws.Range("H44").Value = 14
$ws.Range("I44").Value = 91
$ws.Range("J44").Value = 96.318700000000007
$ws.Range("H45").Value = 14
$ws.Range("I45").Value = 91
$ws.Range("J45").Value = 95.824200000000005
$ws.Range("H46").Value = 14
$ws.Range("I46").Value = 91
$ws.Range("J46").Value = 95.868099999999998
$ws.Range("H47").Value = 14
$ws.Range("I47").Value = 91
$ws.Range("J47").Value = 94.351600000000005
$ws.Range("H48").Value = 14
$ws.Range("I48").Value = 89
$ws.Range("J48").Value = 94.406599999999997
$ws.Range("H49").Value = 14
$ws.Range("I49").Value = 96
$ws.Range("J49").Value = 97.384600000000006
$ws.Range("H50").Value = 14
$ws.Range("I50").Value = 90
$ws.Range("J50").Value = 94.175799999999995
$ws.Range("H51").Value = 14
$ws.Range("I51").Value = 89
$ws.Range("J51").Value = 93.802199999999999

# --- Adjust column widths to (approximately) match the widened/auto-fit columns ---
$ws.Columns.Item(1).ColumnWidth = 2.5924479166666665
$ws.Columns.Item(2).ColumnWidth = 18.022135416666668
$ws.Columns.Item(3).ColumnWidth = 31.022135416666668
$ws.Columns.Item(4).ColumnWidth = 15.307291666666666
$ws.Columns.Item(5).ColumnWidth = 17.736979166666668
$ws.Columns.Item(6).ColumnWidth = 30.736979166666668
$ws.Columns.Item(7).ColumnWidth = 15.022135416666666
$ws.Columns.Item(8).ColumnWidth = 20.451822916666668
$ws.Columns.Item(9).ColumnWidth = 32.736979166666664
$ws.Columns.Item(10).ColumnWidth = 18.736979166666668

# --- Update the view: scroll so column D is at top-left and select I6 ---
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("I6").Select()

